$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '60.514.54'
$ws.Range('E2').Value = '  +2.22%  '
$ws.Range('D3').Value = '2.589.00'
$ws.Range('E3').Value = '  +2.09%  '
$ws.Range('E4').Value = '  -0.25%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '507.19'
$ws.Range('E5').Value = '  +0.67%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '153.83'
$ws.Range('E6').Value = '  -1.25%  '
$ws.Range('E7').Value = '  +0.60%  '
$ws.Range('D9').Value = '2.593.84'
$ws.Range('E9').Value = '  +0.75%  '
$ws.Range('E10').Value = '  +6.74%  '
$ws.Range('E11').Value = '  +0.91%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.347'
$ws.Range('E12').Value = '  +2.40%  '
$ws.Range('E13').Value = '  +0.89%  '
$ws.Range('D14').Value = '3.041.18'
$ws.Range('E14').Value = '  +2.29%  '
$ws.Range('D15').Value = '60.473.20'
$ws.Range('E15').Value = '  +2.33%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '21.49'
$ws.Range('E16').Value = '  -0.97%  '
$ws.Range('E17').Value = '  +2.10%  '
$ws.Range('D18').Value = '2.590.80'
$ws.Range('E18').Value = '  +0.89%  '
$ws.Range('E19').Value = '  +0.95%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '345.76'
$ws.Range('E20').Value = '  +3.75%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '10.44'
$ws.Range('E21').Value = '  +1.36%  '
$ws.Range('E22').Value = '  +1.48%  '
$ws.Range('E23').Value = '  -0.60%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '60.00'
$ws.Range('E24').Value = '  +0.63%  '
$ws.Range('E25').Value = '  +1.26%  '
$ws.Range('E26').Value = '  +0.24%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.998'
$ws.Range('E27').Value = '  +0.34%  '
$ws.Range('D28').Value = '0.0₃0846'
$ws.Range('E28').Value = '  +2.53%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.33'
$ws.Range('E29').Value = '  -0.75%  '
$ws.Range('E30').Value = '  +0.31%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '19.34'
$ws.Range('E31').Value = '  +0.34%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '153.50'
$ws.Range('E32').Value = '  -2.51%  '
$ws.Range('E33').Value = '  -0.76%  '
$ws.Range('E34').Value = '  +3.81%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '3.98'
$ws.Range('E35').Value = '  +2.48%  '
$ws.Range('E36').Value = '  +0.20%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.862'
$ws.Range('E37').Value = '  +13.28%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.851'
$ws.Range('E38').Value = '  +0.64%  '
$ws.Range('B39').Value = 'Filecoin'
$ws.Range('C39').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.76'
$ws.Range('E39').Value = '  +0.63%  '
$ws.Range('B40').Value = 'Stacks'
$ws.Range('C40').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.47'
$ws.Range('E40').Value = '  +1.76%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '35.90'
$ws.Range('E41').Value = '  +2.31%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '296.79'
$ws.Range('E42').Value = '  +2.31%  '
$ws.Range('E43').Value = '  -1.77%  '
$ws.Range('E44').Value = '  -1.50%  '
$ws.Range('B45').Value = 'Hedera'
$ws.Range('C45').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0557'
$ws.Range('E45').Value = '  -1.00%  '
$ws.Range('B46').Value = 'FirstDigitalUSD'
$ws.Range('C46').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.996'
$ws.Range('E46').Value = '  +0.68%  '
$ws.Range('E47').Value = '  +3.86%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '4.85'
$ws.Range('E48').Value = '  +0.40%  '
$ws.Range('E49').Value = '  -0.97%  '
$ws.Range('E50').Value = '  +0.70%  '
$ws.Range('D51').Value = '2.002.99'
$ws.Range('E51').Value = '  +0.36%  '
